# update the start page
#
# 1) Rebrand the title on slide 1 ("EduMate ..." -> "teachMate ...")
# 2) Rebrand the title on slide 3 ("... EduMate" -> "... teachMate")
# 3) Refresh the cached "automatic date" placeholder text (slide master +
#    every slide layout) from 10/17/2025 to 12/9/2025.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 title
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$title1 = $slide1.Shapes.Item(1).TextFrame.TextRange

$title1.Text = "teach"
[void]$title1.InsertAfter("Mate " + [char]0x2013 + " AI-")
[void]$title1.InsertAfter("drivet")
[void]$title1.InsertAfter(" system f" + [char]0x00F6 + "r ")
[void]$title1.InsertAfter("bed" + [char]0x00F6 + "mning")
[void]$title1.InsertAfter(" och feedback")

# ---------------------------------------------------------------------
# 2) Slide 3 title
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(1).TextFrame.TextRange

$title3.Text = "{0}{1}{2}{3}" -f "2", [char]0xFE0F, [char]0x20E3, " "
[void]$title3.InsertAfter("L" + [char]0x00F6 + "sningen")
[void]$title3.InsertAfter(" " + [char]0x2013 + " ")
[void]$title3.InsertAfter("teach")
[void]$title3.InsertAfter("Mate")

# ---------------------------------------------------------------------
# 3) Refresh the "automatically updating" date placeholder text
#    (ppPlaceholderDate = 16) everywhere it appears: the slide master
#    and every slide layout.
# ---------------------------------------------------------------------
$newDate = "12/9/2025"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j)
}
